# Updated cryptos list - refresh Price (col D) and Volume(1h) (col E)
# values on the "cryptos" sheet to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A few Price cells keep a significant trailing zero (e.g. "0.06680")
# that Excel would otherwise drop if the string were auto-detected as a
# number, so force those specific cells to Text before writing them.
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"

$ws.Range("D2").Value = "27.323.75"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "1.713.49"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "224.85"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").Value = "0.5298"
$ws.Range("E6").Value = "  -0.78%  "
$ws.Range("D8").Value = "0.06680"
$ws.Range("E8").Value = "  +1.35%  "
$ws.Range("D9").Value = "0.2650"
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("D10").Value = "20.88"
$ws.Range("E10").Value = "  -2.82%  "
$ws.Range("D11").Value = "0.07703"
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("D12").Value = "4.481"
$ws.Range("E12").Value = "  -2.55%  "
$ws.Range("D13").Value = "1.952.37"
$ws.Range("E13").Value = "  -0.40%  "
$ws.Range("D14").Value = "1.715.03"
$ws.Range("E14").Value = "  -0.64%  "
$ws.Range("D15").Value = "0.5785"
$ws.Range("E15").Value = "  -0.01%  "
$ws.Range("D16").Value = "0.0₅8181"
$ws.Range("E16").Value = "  -1.14%  "
$ws.Range("D17").Value = "67.69"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").Value = "27.360.37"
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("D19").Value = "219.75"
$ws.Range("E19").Value = "  +1.18%  "
$ws.Range("D20").Value = "1.006"
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("D21").Value = "4.642"
$ws.Range("E21").Value = "  -1.47%  "
$ws.Range("D22").Value = "10.42"
$ws.Range("E22").Value = "  -1.28%  "
$ws.Range("D23").Value = "6.021"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").Value = "145.50"
$ws.Range("E25").Value = "  +0.94%  "
$ws.Range("D26").Value = "1.707"
$ws.Range("E26").Value = "  -1.57%  "
$ws.Range("D27").Value = "0.1207"
$ws.Range("E27").Value = "  -1.77%  "
$ws.Range("D28").Value = "7.238"
$ws.Range("E28").Value = "  -1.05%  "
$ws.Range("E29").Value = "  -1.44%  "
$ws.Range("D30").Value = "0.05378"
$ws.Range("E30").Value = "  -1.55%  "
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("D32").Value = "3.478"
$ws.Range("E32").Value = "  -1.86%  "
$ws.Range("D33").Value = "3.391"
$ws.Range("E33").Value = "  -1.12%  "
$ws.Range("D34").Value = "1.636"
$ws.Range("E34").Value = "  -1.18%  "
$ws.Range("D35").Value = "2.852"
$ws.Range("E35").Value = "  -0.42%  "
$ws.Range("D36").Value = "0.9515"
$ws.Range("E36").Value = "  -0.55%  "
$ws.Range("D37").Value = "2.398"
$ws.Range("E37").Value = "  -1.22%  "
$ws.Range("D38").Value = "0.5892"
$ws.Range("E38").Value = "  -0.44%  "
$ws.Range("D39").Value = "1.158.75"
$ws.Range("E39").Value = "  +10.60%  "
$ws.Range("D40").Value = "0.01649"
$ws.Range("E40").Value = "  +0.37%  "
$ws.Range("D41").Value = "5.835"
$ws.Range("E41").Value = "  -1.20%  "
$ws.Range("E42").Value = "  +0.23%  "
$ws.Range("D43").Value = "0.8406"
$ws.Range("E43").Value = "  -0.63%  "
$ws.Range("D44").Value = "100.95"
$ws.Range("E44").Value = "  -0.21%  "
$ws.Range("D45").Value = "1.859.07"
$ws.Range("E45").Value = "  -0.44%  "
$ws.Range("D46").Value = "0.0₈119"
$ws.Range("E46").Value = "  +2.51%  "
$ws.Range("D47").Value = "57.71"
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("D48").Value = "0.4567"
$ws.Range("E48").Value = "  +1.37%  "
$ws.Range("D49").Value = "8.155"
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("D51").Value = "0.05194"
$ws.Range("E51").Value = "  -1.06%  "
